$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (command name) / column E (comments) on Sheet1 is a sorted
# command reference list. "/reportaa" belongs alphabetically right
# between the existing "/quickburns" and "/restock" rows, which sit on
# row 59 before this edit, so insert a fresh row there and push the
# rest of the table down.
$ws.Rows("59:59").Insert()

$ws.Range("A59").Value = "/reportaa"

# Every other command row carries a (possibly blank) comments cell in
# column E; touch it so it is materialized like its neighbours even
# though there's no comment text yet for this new command.
$ws.Range("A59").Style = "Normal"
$ws.Range("E59").Style = "Normal"

Write-Host "Inserted /reportaa at row 59"
